{"js": "// The template paragraph contains:\n//   \"... \u03c4\u03b7\u03c2 \u03c7\u03ce\u03c1\u03b1\u03c2 ${country} \u03b1\u03c0\u03cc ${transport_start} \u03ad\u03c9\u03c2 \u03ba\u03b1\u03b9 ${transport_end} ...\"\n// and the edit wraps the `${country}` placeholder in Greek guillemets, turning it\n// into `\u00ab${country}\u00bb`, matching the style already used for the other placeholders\n// in the same sentence (e.g. \u00ab${teachers}\u00bb, \u00ab${students}\u00bb, \u00ab${school}\u00bb, ...).\n//\n// We locate the unique \" ${country}\" run of text (the leading space distinguishes\n// it from the *other* \"${country}\" placeholder elsewhere in the document, which is\n// already wrapped in guillemets and therefore has no leading space before \"${\").\n// We then insert \"\u00ab\" immediately before \"${\" and, separately, locate the unique\n// \"} \u03b1\u03c0\u03cc \" text and insert \"\u00bb\" immediately after the \"}\".\n\nconst body = context.document.body;\n\n// --- Step 1: insert the opening guillemet \"\u00ab\" right before \"${\" in \" ${country}\" ---\nconst openMatches = body.search(\" ${country}\", { matchCase: true, matchWholeWord: false });\nopenMatches.load(\"text\");\nawait context.sync();\n\nif (openMatches.items.length !== 1) {\n  throw new Error(\"Expected exactly one ' ${country}' match, found \" + openMatches.items.length);\n}\n\nconst openOuter = openMatches.items[0];\nconst dollarBrace = openOuter.search(\"${\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\n\nif (dollarBrace.items.length !== 1) {\n  throw new Error(\"Expected exactly one '${' inside the match, found \" + dollarBrace.items.length);\n}\n\ndollarBrace.items[0].insertText(\"\u00ab\", Word.InsertLocation.before);\nawait context.sync();\n\n// --- Step 2: insert the closing guillemet \"\u00bb\" right after the \"}\" in \"} \u03b1\u03c0\u03cc \" ---\nconst closeMatches = body.search(\"} \u03b1\u03c0\u03cc \", { matchCase: true, matchWholeWord: false });\ncloseMatches.load(\"text\");\nawait context.sync();\n\nif (closeMatches.items.length !== 1) {\n  throw new Error(\"Expected exactly one '} \u03b1\u03c0\u03cc ' match, found \" + closeMatches.items.length);\n}\n\nconst closeOuter = closeMatches.items[0];\nconst closeBrace = closeOuter.search(\"}\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\n\nif (closeBrace.items.length !== 1) {\n  throw new Error(\"Expected exactly one '}' inside the match, found \" + closeBrace.items.length);\n}\n\ncloseBrace.items[0].insertText(\"\u00bb\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The template paragraph contains:\n#   \"... \u03c4\u03b7\u03c2 \u03c7\u03ce\u03c1\u03b1\u03c2 ${country} \u03b1\u03c0\u03cc ${transport_start} \u03ad\u03c9\u03c2 \u03ba\u03b1\u03b9 ${transport_end} ...\"\n# and the edit wraps the `${country}` placeholder in Greek guillemets, turning it\n# into `\u00ab${country}\u00bb`, matching the style already used for the other placeholders\n# in the same sentence (e.g. \u00ab${teachers}\u00bb, \u00ab${students}\u00bb, \u00ab${school}\u00bb, ...).\n#\n# We locate the unique \" ${country}\" text (the leading space distinguishes it from\n# the *other* \"${country}\" placeholder elsewhere in the document, which is already\n# wrapped in guillemets and therefore has no leading space before \"${\"). We then\n# insert \"\u00ab\" immediately before \"${\" and, separately, locate the unique \"} \u03b1\u03c0\u03cc \"\n# text and insert \"\u00bb\" immediately after the \"}\".\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the opening guillemet \"\u00ab\" right before \"${\" in \" ${country}\" ---\n$outer = $d.Content\n$outer.Start = 0\n$outer.End = $d.Content.End\nif (-not $outer.Find.Execute(\" `${country}\")) {\n    throw \"Could not find ' `${country}' in the document\"\n}\n\n$dollarBrace = $outer.Duplicate\nif (-not $dollarBrace.Find.Execute(\"`${\")) {\n    throw \"Could not find '`${' inside the matched range\"\n}\n\n$dollarBrace.InsertBefore(\"\u00ab\")\n\n# --- Step 2: insert the closing guillemet \"\u00bb\" right after the \"}\" in \"} \u03b1\u03c0\u03cc \" ---\n$outer2 = $d.Content\n$outer2.Start = 0\n$outer2.End = $d.Content.End\nif (-not $outer2.Find.Execute(\"} \u03b1\u03c0\u03cc \")) {\n    throw \"Could not find '} \u03b1\u03c0\u03cc ' in the document\"\n}\n\n$closeBrace = $outer2.Duplicate\n$closeBrace.Collapse(1)      # wdCollapseStart\n$closeBrace.MoveEnd(1, 1)    # wdCharacter -> grow to cover just the \"}\"\n$closeBrace.InsertAfter(\"\u00bb\")\n"}
